$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting Peptide Sequence/Modification
# Type/Modification Position/MHC Name (and their data) one column to the
# right.
$ws.Columns.Item(1).Insert()

# Give the new A3 the same formatting (fill + border) as the rest of row 3,
# which the insert did not extend into the brand-new column.
$ws.Range("B3").Copy()
$ws.Range("A3").PasteSpecial(-4122)

# Populate the new column A with the MHC Molecule header + data.
$ws.Range("A1").Value = "MHC Molecule"
$ws.Range("A2").Value = "HLA-A*02:01"
$ws.Range("A3").Value = "HLA-A*02:01"

# The old "MHC Name" column (now shifted to E) is no longer needed since its
# values have been consolidated into the new column A; remove it.
$ws.Columns.Item(5).Delete()

# Update the active selection to match the new cursor position.
$ws.Range("A2").Select()
